$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows 2-11 map to employee_id, employee_name, department, absence_reason,
# absence_duration, absence_date (serial number), salary
$data = @(
    @{ Row=2;  A=39013; B="Sr. Murilo Fonseca";     C="Operacoes";              D="Problemas pessoais"; E=4; F=45087; G=3810.89 }
    @{ Row=3;  A=47813; B="Dr. Rhavi Duarte";        C="Operacoes";              D="Viagem de negocios";  E=3; F=45106; G=9532.799999999999 }
    @{ Row=4;  A=65225; B="Gustavo Henrique Macedo"; C="Vendas";                 D="Problemas pessoais"; E=8; F=45083; G=2397.96 }
    @{ Row=5;  A=39268; B="Cauã Correia";            C="Engenharia";             D="Viagem de negocios";  E=5; F=45087; G=5645 }
    @{ Row=6;  A=67128; B="Manuella Ferreira";       C="Vendas";                 D="Doenca";               E=5; F=45087; G=7082 }
    @{ Row=7;  A=73654; B="Melissa Moreira";         C="Atendimento ao Cliente"; D="Problemas pessoais"; E=6; F=45102; G=3976.85 }
    @{ Row=8;  A=49719; B="Mirella Fogaça";          C="P&D";                    D="Problemas pessoais"; E=6; F=45090; G=5897.83 }
    @{ Row=9;  A=67654; B="Eduarda Cirino";          C="Atendimento ao Cliente"; D="Consulta medica";     E=1; F=45088; G=5596.96 }
    @{ Row=10; A=64407; B="Antônio Moreira";         C="Juridico";               D="Outros";               E=6; F=45093; G=3468.37 }
    @{ Row=11; A=66645; B="Ana Lívia Silveira";      C="TI";                     D="Doenca";               E=2; F=45083; G=2603.51 }
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 1).Value = $item.A
    $ws.Cells.Item($r, 2).Value = $item.B
    $ws.Cells.Item($r, 3).Value = $item.C
    $ws.Cells.Item($r, 4).Value = $item.D
    $ws.Cells.Item($r, 5).Value = $item.E
    $ws.Cells.Item($r, 6).Value = $item.F
    $ws.Cells.Item($r, 7).Value = $item.G
}

$wb.Save()
